$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$f = '=FILTER(LAMBDA(Data, Data)(B3:B24), LAMBDA(Data, Data)(B3:B24) = "A")'
$ws.Range("P2:T10").FormulaArray = $f
